# v2p14. Compatible with MF-Swift v2212, updated hardpoints.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sedan_HambaLG_f")
$ws2 = $wb.Worksheets.Item("Sedan_HambaLG_r")

# Update hardpoint values - Sedan_HambaLG_f (sheet1)
$ws1.Range("F5").Value = -0.0026557142857142869
$ws1.Range("G5").Value = 0.62
$ws1.Range("H5").Value = 0.65

$ws1.Range("F6").Value = 0.055166428571428582
$ws1.Range("G6").Value = 0.85
$ws1.Range("H6").Value = 0.19

# Update hardpoint values - Sedan_HambaLG_r (sheet2)
$ws2.Range("F5").Value = 0.0026557142857142869
$ws2.Range("G5").Value = 0.62
$ws2.Range("H5").Value = 0.65

$ws2.Range("F6").Value = -0.055166428571428582
$ws2.Range("G6").Value = 0.85
$ws2.Range("H6").Value = 0.19

# New number formats applied to the updated cells
$ws1.Range("F5:F6").NumberFormat = "0.000"
$ws1.Range("G5:H6").NumberFormat = "0.00"

$ws2.Range("F5:F6").NumberFormat = "0.000"
$ws2.Range("G5:H6").NumberFormat = "0.00"
